# Apply the commit:
#  - Insert a new row at 43 (pushing old rows 43..60 down to 44..61),
#    and populate it with the content that used to live in row 42
#    (BCIO:036042 "physical performance behaviour", LSR 2 / Intervention
#    content and delivery / External / PS).
#  - Replace row 42's own content with the new BCIO:050432
#    "physical exertion expended on a behaviour" entry (clearing the
#    LSR no. / Ontology section cells and changing Curator to "MS").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a blank row before the current row 43; this shifts the old
#    row 43 (and everything below it) down by one, all formatting intact.
$ws.Rows.Item(43).Insert()

# 2. Fill the newly inserted row 43 with the content that row 42 used
#    to hold before this edit (unstyled, same as the row above it).
$ws.Range("A43").Value = "BCIO:036042"
$ws.Range("B43").Value = "physical performance behaviour"
$ws.Range("C43").Value = "An individual human behaviour that involves maintenance or improvement of flexibility, strength, balance or cardiovascular fitness."
$ws.Range("D43").Value = "individual human behaviour"
$ws.Range("P43").Value = "LSR 2"
$ws.Range("Q43").Value = "Intervention content and delivery"
$ws.Range("S43").Value = "External"
$ws.Range("V43").Value = "PS"

# 3. Overwrite row 42 in place with the new entry.
$ws.Range("A42").Value = "BCIO:050432"
$ws.Range("B42").Value = "physical exertion expended on a behaviour"
$ws.Range("C42").Value = "A behavioural attribute that is the level of musculoskeletal work expended on the behaviour to be enacted."
$ws.Range("D42").Value = "behavioural attribute"
$ws.Range("P42").Value = ""
$ws.Range("Q42").Value = ""
$ws.Range("V42").Value = "MS"
